# Scrum.xlsx update: "Updated scrum and some of the reviews"
$wb = $excel.ActiveWorkbook

# --- Burndown Chart sheet: log 2 hours against "Clean code" on the last day (G column) ---
$wsBurn = $wb.Worksheets.Item("Burndown Chart")
$wsBurn.Range("G14").Value = 2

# --- "01.05" sheet: mark task in row 7 (Observer design pattern doc) as Done ---
$ws0105 = $wb.Worksheets.Item("01.05")
$ws0105.Range("F7").Value = "Done"

# --- View state: move selection on "01.05" off of D10, leave that sheet not active ---
$ws0105.Range("C12").Select()

# --- View state: make "Burndown Chart" the active sheet/tab, with cell O27 selected ---
$wsBurn.Activate()
$wsBurn.Range("O27").Select()
